$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G4").Value = 3.2
$ws.Range("I4").Value = 2.3
$ws.Range("L4").Value = 3.1
$ws.Range("AK4").Value = 21
$ws.Range("AX4").Value = 13

$ws.Range("L5").Value = 4.33
$ws.Range("X5").Value = 9
$ws.Range("AE5").Value = 19
$ws.Range("AO5").Value = 12

$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 9
$ws.Range("Q6").Value = 2.25
$ws.Range("R6").Value = 1.62
